$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set E22 and E24 grading values to 10
$ws.Range("E22").Value = 10
$ws.Range("E24").Value = 10

# Update the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("F22").Select()
